$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-05 10:32:28'
$ws.Range('E3').Value = '2026-02-05 10:32:31'
$ws.Range('E4').Value = '2026-02-05 10:32:33'
$ws.Range('J4').Value = '993.3 hPa'
$ws.Range('K4').Value = '0.4 MJ/m2'
$ws.Range('M4').Value = '10.8 °C 8:34 TU'
$ws.Range('O4').Value = '8.4 °C'
$ws.Range('E5').Value = '2026-02-05 10:32:36'
$ws.Range('J5').Value = '993.5 hPa'
$ws.Range('K5').Value = '0.6 MJ/m2'
$ws.Range('L5').Value = '18.0 km/h - 262º 9:17 TU'
$ws.Range('M5').Value = '9.0 °C 9:19 TU'
$ws.Range('O5').Value = '5.3 °C'
$ws.Range('E6').Value = '2026-02-05 10:32:38'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '68%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H6').PasteSpecial(-4122) | Out-Null
$ws.Range('I6').Value = '0.4 mm'
$ws.Range('J6').Value = '994.2 hPa'
$ws.Range('K6').Value = '0.3 MJ/m2'
$ws.Range('L6').Value = '16.6 km/h - 232º 8:14 TU'
$ws.Range('E7').Value = '2026-02-05 10:32:40'
$ws.Range('E8').Value = '2026-02-05 10:32:43'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '96%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H8').PasteSpecial(-4122) | Out-Null
$ws.Range('K8').Value = '0.5 MJ/m2'
$ws.Range('L8').Value = '9.4 km/h - 118º 9:08 TU'
$ws.Range('M8').Value = '10.0 °C 9:29 TU'
$ws.Range('O8').Value = '4.1 °C'
$ws.Range('E9').Value = '2026-02-05 10:32:45'
$ws.Range('I9').Value = '0.2 mm'
$ws.Range('M9').Value = '3.0 °C 9:54 TU'
$ws.Range('O9').Value = '0.2 °C'
$ws.Range('E10').Value = '2026-02-05 10:32:48'
$ws.Range('M10').Value = '5.3 °C 9:28 TU'
$ws.Range('O10').Value = '2.9 °C'
$ws.Range('E11').Value = '2026-02-05 10:32:50'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '95%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H11').PasteSpecial(-4122) | Out-Null
$ws.Range('J11').Value = '998.5 hPa'
$ws.Range('K11').Value = '0.7 MJ/m2'
$ws.Range('M11').Value = '1.1 °C 9:29 TU'
$ws.Range('O11').Value = '-1.3 °C'
$ws.Range('E12').Value = '2026-02-05 10:32:52'
$ws.Range('E13').Value = '2026-02-05 10:32:55'
$ws.Range('M13').Value = '8.4 °C 9:28 TU'
$ws.Range('O13').Value = '3.9 °C'
$ws.Range('E14').Value = '2026-02-05 10:32:57'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '62%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H14').PasteSpecial(-4122) | Out-Null
$ws.Range('M14').Value = '-1.5 °C 3:08 TU'
$ws.Range('O14').Value = '-2.5 °C'
$ws.Range('E15').Value = '2026-02-05 10:33:00'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '99%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H15').PasteSpecial(-4122) | Out-Null
$ws.Range('J15').Value = '993.9 hPa'
$ws.Range('K15').Value = '0.4 MJ/m2'
$ws.Range('L15').Value = '10.1 km/h - 158º 8:30 TU'
$ws.Range('M15').Value = '8.4 °C 8:59 TU'
$ws.Range('O15').Value = '2.0 °C'
$ws.Range('E16').Value = '2026-02-05 10:33:02'
$ws.Range('E17').Value = '2026-02-05 10:33:04'
$ws.Range('J17').Value = '997.7 hPa'
$ws.Range('K17').Value = '0.4 MJ/m2'
$ws.Range('M17').Value = '2.2 °C 9:04 TU'
$ws.Range('O17').Value = '0.2 °C'
$ws.Range('E18').Value = '2026-02-05 10:33:07'
$ws.Range('E19').Value = '2026-02-05 10:33:09'
$ws.Range('I19').Value = '3.7 mm'
$ws.Range('J19').Value = '995.3 hPa'
$ws.Range('K19').Value = '0.1 MJ/m2'
$ws.Range('L19').Value = '20.5 km/h - 91º 9:28 TU'
$ws.Range('M19').Value = '6.7 °C 9:14 TU'
$ws.Range('O19').Value = '4.8 °C'
$ws.Range('E20').Value = '2026-02-05 10:33:12'
$ws.Range('E21').Value = '2026-02-05 10:33:14'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '92%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H21').PasteSpecial(-4122) | Out-Null
$ws.Range('J21').Value = '994.6 hPa'
$ws.Range('K21').Value = '0.8 MJ/m2'
$ws.Range('L21').Value = '29.2 km/h - 172º 9:47 TU'
$ws.Range('M21').Value = '11.0 °C 9:41 TU'
$ws.Range('O21').Value = '1.9 °C'
$ws.Range('E22').Value = '2026-02-05 10:33:17'
$ws.Range('E23').Value = '2026-02-05 10:33:19'
$ws.Range('E24').Value = '2026-02-05 10:33:22'
$ws.Range('J24').Value = '992.8 hPa'
$ws.Range('K24').Value = '0.6 MJ/m2'
$ws.Range('L24').Value = '40.0 km/h - 202º 9:29 TU'
$ws.Range('E25').Value = '2026-02-05 10:33:24'
$ws.Range('J25').Value = '997.2 hPa'
$ws.Range('K25').Value = '0.8 MJ/m2'
$ws.Range('M25').Value = '1.3 °C 9:37 TU'
$ws.Range('O25').Value = '-0.4 °C'
$ws.Range('E26').Value = '2026-02-05 10:33:27'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '66%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H26').PasteSpecial(-4122) | Out-Null
$ws.Range('K26').Value = '0.7 MJ/m2'
$ws.Range('L26').Value = '22.3 km/h - 274º 9:28 TU'
$ws.Range('M26').Value = '1.8 °C 8:28 TU'
$ws.Range('O26').Value = '-2.6 °C'
$ws.Range('E27').Value = '2026-02-05 10:33:29'
$ws.Range('J27').Value = '993.5 hPa'
$ws.Range('K27').Value = '0.6 MJ/m2'
$ws.Range('M27').Value = '7.4 °C 9:28 TU'
$ws.Range('O27').Value = '4.4 °C'
$ws.Range('E28').Value = '2026-02-05 10:33:32'
$ws.Range('J28').Value = '997.4 hPa'
$ws.Range('L28').Value = '32.4 km/h - 236º 8:40 TU'
$ws.Range('M28').Value = '3.2 °C 9:29 TU'
$ws.Range('O28').Value = '-0.9 °C'
$ws.Range('E29').Value = '2026-02-05 10:33:34'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '91%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H29').PasteSpecial(-4122) | Out-Null
$ws.Range('K29').Value = '0.4 MJ/m2'
$ws.Range('M29').Value = '7.9 °C 9:29 TU'
$ws.Range('O29').Value = '4.9 °C'
$ws.Range('E30').Value = '2026-02-05 10:33:37'
$ws.Range('E31').Value = '2026-02-05 10:33:39'
$ws.Range('E32').Value = '2026-02-05 10:33:41'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '87%'
$ws.Range('H2').Copy() | Out-Null
$ws.Range('H32').PasteSpecial(-4122) | Out-Null
$ws.Range('I32').Value = '0.5 mm'
$ws.Range('J32').Value = '994.2 hPa'
$ws.Range('K32').Value = '0.2 MJ/m2'
$ws.Range('O32').Value = '9.1 °C'
$ws.Range('E33').Value = '2026-02-05 10:33:44'
$ws.Range('M33').Value = '7.8 °C 9:29 TU'
$ws.Range('O33').Value = '4.0 °C'
$ws.Range('E34').Value = '2026-02-05 10:33:46'
$ws.Range('E35').Value = '2026-02-05 10:33:49'
$ws.Range('E36').Value = '2026-02-05 10:33:51'
$ws.Range('I36').Value = '2.0 mm'
$ws.Range('J36').Value = '995.3 hPa'
$ws.Range('K36').Value = '0.2 MJ/m2'
$ws.Range('M36').Value = '7.2 °C 9:29 TU'
$ws.Range('O36').Value = '5.6 °C'
$excel.CutCopyMode = 0
